$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cantines import: replace the terse English enum values with their
# "verbose" French equivalents in row 2 (type_production, type_gestion,
# modèle_économique columns G, H, I).
$ws.Range("G2").Value = "Restaurant avec cuisine sur place"
$ws.Range("H2").Value = "Concédée"
$ws.Range("I2").Value = "Public"

# type_production (G2) keeps its text format but now renders in
# Times New Roman instead of the sheet's default Arial.
$ws.Range("G2").Font.Name = "Times New Roman"

# type_gestion / modèle_économique (H2:I2) switch from the text number
# format to General and now wrap their (longer) text.
$ws.Range("H2:I2").NumberFormat = "General"
$ws.Range("H2:I2").WrapText = $true

# The active selection moves from F2 to I2.
$ws.Range("I2").Select() | Out-Null
